$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the SECOND occurrence of the "Investigacion de nuevas
# tecnologias..." paragraph (the list has two identical copies) gets
# replaced with a paragraph about legal/regulatory research.
# ---------------------------------------------------------------------------
$oldText1 = "Investigación de nuevas tecnologías para incorporar en el aplicativo, como la inteligencia artificial para ofrecer sugerencias personalizadas de cuidado de mascotas."
$newText1 = "Investigación sobre los requisitos legales y regulatorios que podrían afectar el desarrollo y lanzamiento del aplicativo, como las regulaciones de privacidad y seguridad de datos."

# Locate the first occurrence (left untouched) so we can search *after* it
# for the second occurrence, which is the one the diff changes.
$firstHit = $d.Content
$firstHit.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$afterFirst = $d.Range($firstHit.End, $d.Content.End)
$afterFirst.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)

# ---------------------------------------------------------------------------
# Change 2: split the "Evaluar nuevas y mejores formas..." run into three
# runs so the word "app" is wrapped with proofErr (gramStart/gramEnd)
# markers, matching what Word's grammar checker inserts.
# ---------------------------------------------------------------------------
$oldText2 = "Evaluar nuevas y mejores formas para mejorar el diseño de la app con el fin de tener una mejor experiencia de usuario y eficiente forma de representar el contenido."

$para2 = $d.Content
$para2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Remove the original single run's text ...
$para2.Delete()

# ... then rebuild it as three runs with proofErr markers bracketing "app".
$collapsed = $d.Range($para2.Start, $para2.Start)
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr>'
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Evaluar nuevas y mejores formas para mejorar el diseño de la </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPr + '<w:t>app</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> con el fin de tener una mejor experiencia de usuario y eficiente forma de representar el contenido.</w:t></w:r>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$collapsed.InsertXML($xmlFrag)
